$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July")

$ws.Range("B2").Value = 1697
$ws.Range("C2").Value = 1086
$ws.Range("D2").Value = 611
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "1.56 : 1"

$ws.Range("B3").Value = 542
$ws.Range("C3").Value = 456
$ws.Range("D3").Value = 86
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "1.19 : 1"

$ws.Range("B4").Value = 1342
$ws.Range("C4").Value = 1196
$ws.Range("D4").Value = 146
$ws.Range("E4").Value = "We borrowerd more than we lent"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "1.12 : 1"

$ws.Range("B5").Value = 90
$ws.Range("C5").Value = 138
$ws.Range("D5").Value = -48
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.65 : 1"

$ws.Range("B6").Value = 1085
$ws.Range("C6").Value = 1259
$ws.Range("D6").Value = -174
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.86 : 1"

$ws.Range("B7").Value = 163
$ws.Range("C7").Value = 242
$ws.Range("D7").Value = -79
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "We lent more than we borrowed"
$ws.Range("G7").Value = "0.67 : 1"

$ws.Range("B8").Value = 128
$ws.Range("C8").Value = 203
$ws.Range("D8").Value = -75
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.63 : 1"

$ws.Range("B9").Value = 54
$ws.Range("C9").Value = 53
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "We borrowerd more than we lent"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = "1.02 : 1"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 36
$ws.Range("D10").Value = -36
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.00 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""

$ws.Range("B12").Value = 22
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = "We borrowerd more than we lent"
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = "2.00 : 1"

$ws.Range("B13").Value = 139
$ws.Range("C13").Value = 68
$ws.Range("D13").Value = 71
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = "2.04 : 1"

$ws.Range("B14").Value = 129
$ws.Range("C14").Value = 309
$ws.Range("D14").Value = -180
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.42 : 1"

$ws.Range("B15").Value = 143
$ws.Range("C15").Value = 125
$ws.Range("D15").Value = 18
$ws.Range("E15").Value = "We borrowerd more than we lent"
$ws.Range("F15").Value = ""
$ws.Range("G15").Value = "1.14 : 1"

$ws.Range("B16").Value = 45
$ws.Range("C16").Value = 172
$ws.Range("D16").Value = -127
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.26 : 1"

$ws.Range("B17").Value = 534
$ws.Range("C17").Value = 423
$ws.Range("D17").Value = 111
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = "1.26 : 1"

$ws.Range("B18").Value = 53
$ws.Range("C18").Value = 84
$ws.Range("D18").Value = -31
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.63 : 1"

$ws.Range("B19").Value = 586
$ws.Range("C19").Value = 419
$ws.Range("D19").Value = 167
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = "1.40 : 1"

$ws.Range("B20").Value = 6
$ws.Range("C20").Value = 62
$ws.Range("D20").Value = -56
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.10 : 1"

$ws.Range("B21").Value = 451
$ws.Range("C21").Value = 404
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("F21").Value = ""
$ws.Range("G21").Value = "1.12 : 1"

$ws.Range("B22").Value = 22
$ws.Range("C22").Value = 90
$ws.Range("D22").Value = -68
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.24 : 1"

$ws.Range("B23").Value = 647
$ws.Range("C23").Value = 435
$ws.Range("D23").Value = 212
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("F23").Value = ""
$ws.Range("G23").Value = "1.49 : 1"

$ws.Range("B24").Value = 1711
$ws.Range("C24").Value = 1453
$ws.Range("D24").Value = 258
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("F24").Value = ""
$ws.Range("G24").Value = "1.18 : 1"

$ws.Range("B25").Value = 230
$ws.Range("C25").Value = 437
$ws.Range("D25").Value = -207
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.53 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = ""

$ws.Range("B27").Value = 313
$ws.Range("C27").Value = 160
$ws.Range("D27").Value = 153
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("F27").Value = ""
$ws.Range("G27").Value = "1.96 : 1"

$ws.Range("B28").Value = 55
$ws.Range("C28").Value = 69
$ws.Range("D28").Value = -14
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = "We lent more than we borrowed"
$ws.Range("G28").Value = "0.80 : 1"

$ws.Range("B29").Value = 503
$ws.Range("C29").Value = 366
$ws.Range("D29").Value = 137
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("F29").Value = ""
$ws.Range("G29").Value = "1.37 : 1"

$ws.Range("B30").Value = 91
$ws.Range("C30").Value = 37
$ws.Range("D30").Value = 54
$ws.Range("E30").Value = "We borrowerd more than we lent"
$ws.Range("F30").Value = ""
$ws.Range("G30").Value = "2.46 : 1"

$ws.Range("B31").Value = 42
$ws.Range("C31").Value = 310
$ws.Range("D31").Value = -268
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.14 : 1"

$ws.Range("B32").Value = 444
$ws.Range("C32").Value = 467
$ws.Range("D32").Value = -23
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.95 : 1"

$ws.Range("B33").Value = 351
$ws.Range("C33").Value = 506
$ws.Range("D33").Value = -155
$ws.Range("E33").Value = ""
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.69 : 1"

$ws.Range("B34").Value = 166
$ws.Range("C34").Value = 84
$ws.Range("D34").Value = 82
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("F34").Value = ""
$ws.Range("G34").Value = "1.98 : 1"

$ws.Range("B35").Value = 910
$ws.Range("C35").Value = 887
$ws.Range("D35").Value = 23
$ws.Range("E35").Value = "We borrowerd more than we lent"
$ws.Range("F35").Value = ""
$ws.Range("G35").Value = "1.03 : 1"

$ws.Range("B36").Value = 246
$ws.Range("C36").Value = 506
$ws.Range("D36").Value = -260
$ws.Range("E36").Value = ""
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.49 : 1"

$ws.Range("B37").Value = 548
$ws.Range("C37").Value = 327
$ws.Range("D37").Value = 221
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("F37").Value = ""
$ws.Range("G37").Value = "1.68 : 1"

$ws.Range("B38").Value = 8
$ws.Range("C38").Value = 141
$ws.Range("D38").Value = -133
$ws.Range("E38").Value = ""
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.06 : 1"

$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 6
$ws.Range("D39").Value = -6
$ws.Range("E39").Value = ""
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.00 : 1"

$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = -1
$ws.Range("E40").Value = ""
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.00 : 1"

$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = -2
$ws.Range("E41").Value = ""
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.00 : 1"

$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = -1
$ws.Range("E42").Value = ""
$ws.Range("F42").Value = "We lent more than we borrowed"
$ws.Range("G42").Value = "0.00 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = ""
$ws.Range("F43").Value = ""
$ws.Range("G43").Value = ""

$ws.Range("B44").Value = 64
$ws.Range("C44").Value = 67
$ws.Range("D44").Value = -3
$ws.Range("E44").Value = ""
$ws.Range("F44").Value = "We lent more than we borrowed"
$ws.Range("G44").Value = "0.96 : 1"

$ws.Range("B45").Value = 75
$ws.Range("C45").Value = 156
$ws.Range("D45").Value = -81
$ws.Range("E45").Value = ""
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.48 : 1"

$ws.Range("B46").Value = 409
$ws.Range("C46").Value = 592
$ws.Range("D46").Value = -183
$ws.Range("E46").Value = ""
$ws.Range("F46").Value = "We lent more than we borrowed"
$ws.Range("G46").Value = "0.69 : 1"

$ws.Range("B47").Value = 732
$ws.Range("C47").Value = 547
$ws.Range("D47").Value = 185
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("F47").Value = ""
$ws.Range("G47").Value = "1.34 : 1"

$ws.Range("B48").Value = 205
$ws.Range("C48").Value = 533
$ws.Range("D48").Value = -328
$ws.Range("E48").Value = ""
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.38 : 1"

$ws.Range("B49").Value = 322
$ws.Range("C49").Value = 218
$ws.Range("D49").Value = 104
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("F49").Value = ""
$ws.Range("G49").Value = "1.48 : 1"

$ws.Range("B50").Value = 771
$ws.Range("C50").Value = 513
$ws.Range("D50").Value = 258
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("F50").Value = ""
$ws.Range("G50").Value = "1.50 : 1"

$ws.Range("B51").Value = 213
$ws.Range("C51").Value = 116
$ws.Range("D51").Value = 97
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("F51").Value = ""
$ws.Range("G51").Value = "1.84 : 1"

$ws.Range("B52").Value = 361
$ws.Range("C52").Value = 608
$ws.Range("D52").Value = -247
$ws.Range("E52").Value = ""
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.59 : 1"

$ws.Range("B53").Value = 194
$ws.Range("C53").Value = 251
$ws.Range("D53").Value = -57
$ws.Range("E53").Value = ""
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.77 : 1"

$ws.Range("B54").Value = 25
$ws.Range("C54").Value = 225
$ws.Range("D54").Value = -200
$ws.Range("E54").Value = ""
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.11 : 1"

$ws.Range("B55").Value = 205
$ws.Range("C55").Value = 215
$ws.Range("D55").Value = -10
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = "We lent more than we borrowed"
$ws.Range("G55").Value = "0.95 : 1"

$ws.Range("G30").Select()